$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.128.34'
$ws.Range("E2").Value = '  -0.68%  '
$ws.Range("D3").Value = '2.446.37'
$ws.Range("E3").Value = '  +0.19%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '582.06'
$ws.Range("E5").Value = '  +1.40%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '142.81'
$ws.Range("E6").Value = '  -1.05%  '
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("E8").Value = '  +0.37%  '
$ws.Range("D9").Value = '2.441.66'
$ws.Range("E9").Value = '  +0.19%  '
$ws.Range("E10").Value = '  +1.99%  '
$ws.Range("E11").Value = '  +2.33%  '
$ws.Range("E12").Value = '  -0.16%  '
$ws.Range("E13").Value = '  -2.67%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.43'
$ws.Range("E14").Value = '  -0.81%  '
$ws.Range("E15").Value = '  -0.05%  '
$ws.Range("D16").Value = '2.878.74'
$ws.Range("E16").Value = '  -0.13%  '
$ws.Range("D17").Value = '61.993.38'
$ws.Range("E17").Value = '  -0.63%  '
$ws.Range("D18").Value = '2.432.87'
$ws.Range("E18").Value = '  -0.33%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.70'
$ws.Range("E19").Value = '  -3.97%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.22'
$ws.Range("E20").Value = '  +0.71%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '326.12'
$ws.Range("E21").Value = '  -0.77%  '
$ws.Range("E22").Value = '  -1.31%  '
$ws.Range("E23").Value = '  +0.01%  '
$ws.Range("E24").Value = '  -6.54%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '65.60'
$ws.Range("E25").Value = '  +0.10%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.13'
$ws.Range("E26").Value = '  -0.77%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '598.00'
$ws.Range("E27").Value = '  -5.84%  '
$ws.Range("D28").Value = '0.0₃0967'
$ws.Range("E28").Value = '  +0.04%  '
$ws.Range("D29").Value = '2.564.70'
$ws.Range("E29").Value = '  +0.23%  '
$ws.Range("E30").Value = '  +1.30%  '
$ws.Range("E31").Value = '  -1.70%  '
$ws.Range("E32").Value = '  -1.46%  '
$ws.Range("E33").Value = '  +1.11%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.134'
$ws.Range("E34").Value = '  -0.96%  '
$ws.Range("E35").Value = '  -2.64%  '
$ws.Range("E36").Value = '  +0.24%  '
$ws.Range("E37").Value = '  -2.07%  '
$ws.Range("E38").Value = '  +0.11%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '153.69'
$ws.Range("E39").Value = '  +5.14%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.29'
$ws.Range("E41").Value = '  +0.60%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '43.25'
$ws.Range("E42").Value = '  +2.32%  '
$ws.Range("E43").Value = '  -1.74%  '
$ws.Range("E44").Value = '  +0.05%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.52'
$ws.Range("E45").Value = '  +0.74%  '
$ws.Range("D46").Value = '0.0₆0275'
$ws.Range("E46").Value = '  +22.36%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '141.73'
$ws.Range("E47").Value = '  -2.65%  '
$ws.Range("E48").Value = '  -2.82%  '
$ws.Range("E49").Value = '  +0.51%  '
$ws.Range("E50").Value = '  -1.29%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '19.84'
$ws.Range("E51").Value = '  +0.74%  '
